$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.119011402130127
$ws.Range("B1").Value = 2.596107959747314
$ws.Range("C1").Value = 2.709514141082764
$ws.Range("D1").Value = 3.107975244522095
$ws.Range("E1").Value = 0.7647616863250732
